# Update the "Förändrad" (column C) date for all existing data rows (2..325)
# from 45203 to 45204.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C325").Value = 45204

# Row 326: beteckning changes, Förändrad date changes, area changes.
$ws.Cells.Item(326, 1).Value = "A 47296-2023"
$ws.Cells.Item(326, 3).Value = 45204
$ws.Cells.Item(326, 7).Value = 0.6

# Row 327: beteckning changes, Förändrad date changes, area changes,
# and the row gains an explicit row height (15pt, matching the other rows).
$ws.Cells.Item(327, 1).Value = "A 47305-2023"
$ws.Cells.Item(327, 3).Value = 45204
$ws.Cells.Item(327, 7).Value = 2
$ws.Rows.Item(327).RowHeight = 15

# New rows 328-332 appended with the same shape/format as the existing rows.
$newRows = @(
    @{ Row = 328; A = "A 47300-2023"; G = 1.8 },
    @{ Row = 329; A = "A 47309-2023"; G = 1.6 },
    @{ Row = 330; A = "A 47307-2023"; G = 2.6 },
    @{ Row = 331; A = "A 47303-2023"; G = 2 },
    @{ Row = 332; A = "A 47310-2023"; G = 0.9 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = 45202
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = 45204
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "JÖNKÖPINGS LÄN"
    $ws.Cells.Item($row, 5).Value = "ANEBY"
    $ws.Cells.Item($row, 6).Value = "Sveaskog"

    $ws.Cells.Item($row, 7).Value = $r.G

    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }

    $ws.Range("R" + $row).WrapText = $true
}

# Rows 328-331 get an explicit 15pt row height (the last row, 332, does not,
# matching the pattern of the previous last row before this edit).
$ws.Range("A328:A331").RowHeight = 15
